$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B27").Value = "c70b2ee0681fdd4002130ecd66ac8f38"
$ws.Range("B193").Value = "d74535d02015b6460edb4478546765b9"
$ws.Range("B280").Value = "c404f70f45265e545815071a601f77d9"
$ws.Range("B297").Value = "c3b26e14fb1d78a2b5bac356595a2526"
$ws.Range("B300").Value = "162f468954607ed65138efbdf594ba88"
$ws.Range("B358").Value = "df6d574de7ef3d9d0de84de64fe3896c"
$ws.Range("B367").Value = "a0000ab18dfbd9c0c767788e147199a0"
$ws.Range("B397").Value = "e899843e8de1d189c9c71a6969ab9f97"
$ws.Range("B416").Value = "07de456c71140bc734272a3527739cdd"
$ws.Range("B422").Value = "f1e38f056238d0018645dbb553ff687e"
$ws.Range("B477").Value = "67996393d700dcfd73f87d83e57729d0"
$ws.Range("B510").Value = "495a926323289283c84bcf78a14431b2"
$ws.Range("B511").Value = "ba8b7906a962dfdff615d1b72c8df159"
$ws.Range("B520").Value = "19d63f95d907d6346857c3d81f5ac1ae"
$ws.Range("B529").Value = "c0011aff164bd146db9a57c424f84090"
$ws.Range("B546").Value = "89f90dca2c5b58d21fc45cdb4cbe515b"
$ws.Range("B564").Value = "1cdae3625b8e712d758e4c08c68d46e2"
$ws.Range("B577").Value = "7abd9dec4bc8df440f329ab58129597a"
$ws.Range("B589").Value = "786809fcbc89ae603c969e3dae04eaa6"
$ws.Range("B770").Value = "817fea29e3dfb322fbfe948fdc078cf0"
$ws.Range("B789").Value = "31d9cdd93eedbe98f1341fbe81922060"
$ws.Range("B803").Value = "1cf281868cee3e059ad0a19345b59263"
$ws.Range("B897").Value = "97b8c77a8451b9f94c2ebdb90798c416"
$ws.Range("B905").Value = "e5bfb00096365feb835b7f6eb5251980"
$ws.Range("B963").Value = "74ce11a521c514d8df914174f6efb73d"
$ws.Range("B967").Value = "1382dc1aa6457e2dfe23d4db3af80247"
